$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - isolator: add Supplier + Price
$ws.Range("C4").Value = "AFW"
$ws.Range("D4").Value = "260 + GST"

# Row 5 - was "Mirror mount (grating)" -> "Mirror mount (grating mount)"
$ws.Range("B5").Value = "Mirror mount (grating mount)"

# Row 6 - stays "Lens tube" (text unchanged, only shared-string index shifts)
$ws.Range("B6").Value = "Lens tube"

# Row 7 - unchanged ("Lens, NA 0.65")

# Row 8 - becomes the new wedged-window design line, with a bold part-number in A8
$ws.Range("A8").Value = "WW11050-C14"
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").Font.Size = 10
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Color = 0x33 + 0x33 * 256 + 0x33 * 65536
$ws.Range("B8").Value = "Wedged N-BK7 Laser Window, 1dia, AR Coated: 1047 - 1064 nm"
$ws.Range("C8").Value = "Thorlabs"
$ws.Range("D8").Value = 120

# Row 9 - diffraction grating line gains a part number, supplier, price
$ws.Range("A9").Value = "GR13-1210"
$ws.Range("B9").Value = "Ruled reflective diffraction grating, 1200/mm, 1um blaze, 12.7x12.7x6mm"
$ws.Range("C9").Value = "Thorlabs"
$ws.Range("D9").Value = 65

# Row 10 - PZT line gains a part number, supplier, price
$ws.Range("A10").Value = "PA4FK"
$ws.Range("B10").Value = "4 x Piezo chip, 150V, 3.6um, bare electrodes"
$ws.Range("C10").Value = "Thorlabs"
$ws.Range("D10").Value = 31

# Row 11 - stays "Fiber connectors"
$ws.Range("B11").Value = "Fiber connectors"

# Row 12 - stays "Electrical connectors"
$ws.Range("B12").Value = "Electrical connectors"
